$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - (Intercept)
$ws.Range("B2").Value = 36838.023985
$ws.Range("D2").Value = 51.418332

# Row 3 - household_group_collapsed
$ws.Range("B3").Value = 12345.922808
$ws.Range("D3").Value = 8.616189
$ws.Range("E3").Value = 0.000249

# Row 4 - Residuals
$ws.Range("B4").Value = 159765.575023
$ws.Range("C4").Value = 223

# Row 5 - SM-Control
$ws.Range("G5").Value = -10.000859
$ws.Range("H5").Value = -21.690692
$ws.Range("I5").Value = 1.688973
$ws.Range("J5").Value = 0.110103

# Row 6 - SM + Traps-Control
$ws.Range("G6").Value = 6.429079
$ws.Range("H6").Value = -6.158395
$ws.Range("I6").Value = 19.016554
$ws.Range("J6").Value = 0.451391

# Row 7 - SM + Traps-SM
$ws.Range("G7").Value = 16.429939
$ws.Range("H7").Value = 6.913677
$ws.Range("I7").Value = 25.9462
$ws.Range("J7").Value = 0.00019
